# Re-colour the deck's (single) design theme to the stock "Office Theme"
# palette. PowerPoint stores theme colours as dk1/lt1/dk2/lt2/accent1-6/
# hlink/folHlink (12 slots) and exposes them through the active slide's
# ThemeColorScheme, which edits the shared theme part used by the slide
# master (and therefore every slide built on it).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# index -> (slot name, new RGB as 0x00BBGGRR for Office theme colours)
$newColors = @(
    0,         # 1  dk1      000000
    16777215,  # 2  lt1      FFFFFF
    6968388,   # 3  dk2      44546A
    15132391,  # 4  lt2      E7E6E6
    13998939,  # 5  accent1  5B9BD5
    3243501,   # 6  accent2  ED7D31
    10855845,  # 7  accent3  A5A5A5
    49407,     # 8  accent4  FFC000
    12874308,  # 9  accent5  4472C4
    4697456,   # 10 accent6  70AD47
    12673797,  # 11 hlink    0563C1
    7491477    # 12 folHlink 954F72
)

for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Item($i).RGB = $newColors[$i - 1]
}
